# Applies the July 07, 2020 09:10:33 AM America/Chicago TZ run results
# to the covid_disparities_output workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Virginia
$ws.Range("B11").Value = 44019
$ws.Range("C11").Value = 66740
$ws.Range("D11").Value = 1881
$ws.Range("E11").Value = 9747
$ws.Range("F11").Value = 418
$ws.Range("G11").Value = 14.6
$ws.Range("H11").Value = 22.22

# Row 17 - Maryland
$ws.Range("B17").Value = 44019
$ws.Range("C17").Value = 70396
$ws.Range("D17").Value = 3140
$ws.Range("E17").Value = 20157
$ws.Range("F17").Value = 1266
$ws.Range("G17").Value = 28.63
$ws.Range("H17").Value = 40.32
$ws.Range("K17").Value = 57574
$ws.Range("L17").Value = 70366

# Row 32 - Washington, DC
$ws.Range("B32").Value = 44019
$ws.Range("C32").Value = 10569
$ws.Range("D32").Value = 561
$ws.Range("E32").Value = 5241
$ws.Range("F32").Value = 418
$ws.Range("G32").Value = 49.59
$ws.Range("H32").Value = 74.51000000000001

# Row 34 - Maine
$ws.Range("B34").Value = 44019
$ws.Range("C34").Value = 3440
$ws.Range("D34").Value = 110
$ws.Range("E34").Value = 802
$ws.Range("G34").Value = 26.42
$ws.Range("K34").Value = 3036
